$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.142.54'
$ws.Range('E2').Value = '  +0.03%  '

# Row 3
$ws.Range('E3').Value = '  -0.69%  '

# Row 4
$ws.Range('D4').Value = "'1.003"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.21%  '

# Row 5
$ws.Range('D5').Value = "'313.04"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.39%  '

# Row 6
$ws.Range('E6').Value = '  +0.19%  '

# Row 7
$ws.Range('D7').Value = "'0.5051"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.60%  '

# Row 8
$ws.Range('D8').Value = "'0.3834"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.62%  '

# Row 9
$ws.Range('D9').Value = "'0.08548"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.49%  '

# Row 10
$ws.Range('E10').Value = '  -1.21%  '

# Row 11
$ws.Range('D11').Value = "'41.76"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.22%  '

# Row 12
$ws.Range('D12').Value = "'6.271"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.01%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.876.30'
$ws.Range('E13').Value = '  -1.60%  '

# Row 14
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').Value = "'20.61"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.00%  '

# Row 15
$ws.Range('D15').Value = "'7.216"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.04%  '

# Row 16
$ws.Range('D16').Value = "'1.003"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.23%  '

# Row 17
$ws.Range('E17').Value = '  -0.96%  '

# Row 18
$ws.Range('D18').Value = "'91.21"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.32%  '

# Row 19
$ws.Range('D19').Value = "'0.06667"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.07%  '

# Row 20
$ws.Range('D20').Value = "'18.11"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.38%  '

# Row 21
$ws.Range('E21').Value = '  +0.27%  '

# Row 22
$ws.Range('D22').Value = "'6.105"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.69%  '

# Row 23
$ws.Range('D23').Value = '28.178.88'
$ws.Range('E23').Value = '  -0.04%  '

# Row 24
$ws.Range('E24').Value = '  -2.35%  '

# Row 25
$ws.Range('D25').Value = "'2.267"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.31%  '

# Row 26
$ws.Range('E26').Value = '  +1.23%  '

# Row 27
$ws.Range('D27').Value = '2.092.85'
$ws.Range('E27').Value = '  -1.61%  '

# Row 28
$ws.Range('E28').Value = '  -0.81%  '

# Row 29
$ws.Range('D29').Value = "'156.41"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.20%  '

# Row 30
$ws.Range('D30').Value = "'126.52"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.32%  '

# Row 31
$ws.Range('E31').Value = '  -0.87%  '

# Row 32
$ws.Range('D32').Value = "'1.054"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.75%  '

# Row 33
$ws.Range('D33').Value = "'5.644"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.61%  '

# Row 34
$ws.Range('D34').Value = "'3.608"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.43%  '

# Row 35
$ws.Range('D35').Value = "'9.707"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.42%  '

# Row 36
$ws.Range('D36').Value = "'0.02458"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.48%  '

# Row 37
$ws.Range('D37').Value = "'0.06554"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.76%  '

# Row 38
$ws.Range('D38').Value = "'1.231"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.46%  '

# Row 39
$ws.Range('D39').Value = "'0.2179"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.24%  '

# Row 40
$ws.Range('D40').Value = "'0.6538"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.05%  '

# Row 41
$ws.Range('D41').Value = "'1.243"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.87%  '

# Row 42
$ws.Range('E42').Value = '  +0.03%  '

# Row 43
$ws.Range('D43').Value = "'4.917"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.25%  '

# Row 44
$ws.Range('D44').Value = "'0.6205"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.60%  '

# Row 45
$ws.Range('D45').Value = "'13.16"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.10%  '

# Row 46
$ws.Range('E46').Value = '  -0.46%  '

# Row 47
$ws.Range('E47').Value = '  -0.23%  '

# Row 48
$ws.Range('D48').Value = "'2.019"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.71%  '

# Row 49
$ws.Range('D49').Value = "'1.219"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.52%  '

# Row 50
$ws.Range('D50').Value = "'121.01"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.02%  '

# Row 51
$ws.Range('D51').Value = "'80.80"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.99%  '
